$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 is an automatic daily update of electricity spot prices.
# Update the date (A2) and all hourly / summary values to the new day's data.

$ws.Range("A2").Value = "2025-12-08"

$ws.Range("B2").Value = 71.34
$ws.Range("C2").Value = 65.53
$ws.Range("D2").Value = 52.42
$ws.Range("E2").Value = 50
$ws.Range("F2").Value = 47.46
$ws.Range("G2").Value = 50.18
$ws.Range("H2").Value = 52.05
$ws.Range("I2").Value = 67.40000000000001
$ws.Range("J2").Value = 64.15000000000001
$ws.Range("K2").Value = 45.21
$ws.Range("L2").Value = 27.82
$ws.Range("M2").Value = 23.02
$ws.Range("N2").Value = 20.92
$ws.Range("O2").Value = 14.12
$ws.Range("P2").Value = 15.54
$ws.Range("Q2").Value = 24.81
$ws.Range("R2").Value = 39.38
$ws.Range("S2").Value = 76.06
$ws.Range("T2").Value = 83.31
$ws.Range("U2").Value = 91.09
$ws.Range("V2").Value = 96.27
$ws.Range("W2").Value = 95.29000000000001
$ws.Range("X2").Value = 81.52
$ws.Range("Y2").Value = 72.88
$ws.Range("Z2").Value = 55.32

# AA2 (Slot_4h_max) unchanged: "20h-24h"
$ws.Range("AB2").Value = 86.48999999999999

# AC2 (Slot_2h_frist) unchanged: "20h-22h"
$ws.Range("AD2").Value = 95.78

# AE2 (Slot_2h_second) unchanged: "18h-20h"
$ws.Range("AF2").Value = 87.2

$ws.Range("AG2").Value = "2h-16h"
